$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (no explicit formatting) used to strip any incidental
# number-format/style that Excel may apply when a numeric-looking string
# is typed into a cell, so the saved style index matches the original.
$defaultStyle = $ws.Range("B2").Style

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $defaultStyle
}

Set-TextValue 'D2' '42.761.82'
Set-TextValue 'E2' '  -7.14%  '
Set-TextValue 'D3' '2.547.28'
Set-TextValue 'E3' '  -2.07%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '296.49'
Set-TextValue 'E5' '  -5.09%  '
Set-TextValue 'D6' '91.48'
Set-TextValue 'E6' '  -7.40%  '
Set-TextValue 'D7' '0.573'
Set-TextValue 'E7' '  -4.38%  '
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'D9' '0.547'
Set-TextValue 'E9' '  -6.07%  '
Set-TextValue 'D10' '35.57'
Set-TextValue 'E10' '  -8.84%  '
Set-TextValue 'D11' '0.0806'
Set-TextValue 'E11' '  -4.10%  '
Set-TextValue 'D12' '7.63'
Set-TextValue 'E12' '  -6.21%  '
Set-TextValue 'D13' '2.933.16'
Set-TextValue 'E13' '  -2.03%  '
Set-TextValue 'E14' '  -0.04%  '
Set-TextValue 'D15' '2.546.25'
Set-TextValue 'E15' '  -4.81%  '
Set-TextValue 'D16' '0.863'
Set-TextValue 'E16' '  -5.82%  '
Set-TextValue 'D17' '14.10'
Set-TextValue 'E17' '  -5.25%  '
Set-TextValue 'D18' '42.748.17'
Set-TextValue 'E18' '  -7.47%  '
Set-TextValue 'B19' 'ShibaInu'
Set-TextValue 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D19' '0.0₃0975'
Set-TextValue 'E19' '  -4.24%  '
Set-TextValue 'B20' 'Uniswap'
Set-TextValue 'C20' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D20' '6.63'
Set-TextValue 'E20' '  -1.38%  '
Set-TextValue 'D21' '12.45'
Set-TextValue 'E21' '  -2.71%  '
Set-TextValue 'D22' '72.51'
Set-TextValue 'E22' '  -0.17%  '
Set-TextValue 'D23' '259.57'
Set-TextValue 'E23' '  -10.70%  '
Set-TextValue 'D24' '2.88'
Set-TextValue 'E24' '  -5.91%  '
Set-TextValue 'D25' '29.51'
Set-TextValue 'E25' '  -1.95%  '
Set-TextValue 'E26' '  -7.63%  '
Set-TextValue 'E27' '  +0.00%  '
Set-TextValue 'D28' '9.97'
Set-TextValue 'E28' '  -7.61%  '
Set-TextValue 'E29' '  -4.48%  '
Set-TextValue 'D30' '36.14'
Set-TextValue 'E30' '  -5.36%  '
Set-TextValue 'D31' '5.89'
Set-TextValue 'E31' '  -5.85%  '
Set-TextValue 'D32' '150.35'
Set-TextValue 'E32' '  -3.46%  '
Set-TextValue 'D33' '2.16'
Set-TextValue 'E33' '  -2.27%  '
Set-TextValue 'D34' '3.39'
Set-TextValue 'E34' '  -5.43%  '
Set-TextValue 'E35' '  -3.60%  '
Set-TextValue 'D36' '0.0792'
Set-TextValue 'E36' '  -5.67%  '
Set-TextValue 'D37' '0.114'
Set-TextValue 'E37' '  -6.67%  '
Set-TextValue 'B38' 'EnergySwap'
Set-TextValue 'C38' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D38' '24.25'
Set-TextValue 'E38' '  +8.69%  '
Set-TextValue 'B39' 'Stellar'
Set-TextValue 'C39' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D39' '0.119'
Set-TextValue 'E39' '  -3.05%  '
Set-TextValue 'D40' '16.19'
Set-TextValue 'E40' '  +2.31%  '
Set-TextValue 'D41' '3.42'
Set-TextValue 'E41' '  -4.84%  '
Set-TextValue 'D42' '0.0309'
Set-TextValue 'E42' '  -6.64%  '
Set-TextValue 'D43' '2.079.36'
Set-TextValue 'E43' '  -0.88%  '
Set-TextValue 'D44' '3.81'
Set-TextValue 'E44' '  -3.65%  '
Set-TextValue 'D45' '0.998'
Set-TextValue 'E45' '  -0.03%  '
Set-TextValue 'D46' '84.70'
Set-TextValue 'E46' '  -13.21%  '
Set-TextValue 'E47' '  +2.94%  '
Set-TextValue 'B48' 'Stacks'
Set-TextValue 'C48' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D48' '1.71'
Set-TextValue 'E48' '  -1.93%  '
Set-TextValue 'B49' 'RocketPoolETH'
Set-TextValue 'C49' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D49' '2.791.52'
Set-TextValue 'E49' '  -2.06%  '
Set-TextValue 'B50' 'FraxShare'
Set-TextValue 'C50' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D50' '8.71'
Set-TextValue 'E50' '  -9.67%  '
Set-TextValue 'D51' '103.65'
Set-TextValue 'E51' '  -4.44%  '
